$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled like the existing H1 header
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I (I0) and J (IF), rows 2-82
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 8
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 9
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 8
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 7
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 5
$ws.Range("I12").Value = 7
$ws.Range("J12").Value = 7
$ws.Range("I13").Value = 7
$ws.Range("J13").Value = 7
$ws.Range("I14").Value = 9
$ws.Range("J14").Value = 9
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 5
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 7
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 5
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 8
$ws.Range("I19").Value = 9
$ws.Range("J19").Value = 9
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 8
$ws.Range("I21").Value = 6
$ws.Range("J21").Value = 7
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 7
$ws.Range("I23").Value = 7
$ws.Range("J23").Value = 7
$ws.Range("I24").Value = 8
$ws.Range("J24").Value = 8
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 5
$ws.Range("I26").Value = 6
$ws.Range("J26").Value = 6
$ws.Range("I27").Value = 6
$ws.Range("J27").Value = 6
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 8
$ws.Range("I29").Value = 9
$ws.Range("J29").Value = 9
$ws.Range("I30").Value = 7
$ws.Range("J30").Value = 8
$ws.Range("I31").Value = 8
$ws.Range("J31").Value = 8
$ws.Range("I32").Value = 8
$ws.Range("J32").Value = 8
$ws.Range("I33").Value = 9
$ws.Range("J33").Value = 9
$ws.Range("I34").Value = 9
$ws.Range("J34").Value = 9
$ws.Range("I35").Value = 8
$ws.Range("J35").Value = 9
$ws.Range("I36").Value = 8
$ws.Range("J36").Value = 8
$ws.Range("I37").Value = 13
$ws.Range("J37").Value = 14
$ws.Range("I38").Value = 8
$ws.Range("J38").Value = 8
$ws.Range("I39").Value = 9
$ws.Range("J39").Value = 9
$ws.Range("I40").Value = 8
$ws.Range("J40").Value = 9
$ws.Range("I41").Value = 9
$ws.Range("J41").Value = 9
$ws.Range("I42").Value = 7
$ws.Range("J42").Value = 8
$ws.Range("I43").Value = 8
$ws.Range("J43").Value = 8
$ws.Range("I44").Value = 7
$ws.Range("J44").Value = 8
$ws.Range("I45").Value = 8
$ws.Range("J45").Value = 8
$ws.Range("I46").Value = 8
$ws.Range("J46").Value = 9
$ws.Range("I47").Value = 9
$ws.Range("J47").Value = 9
$ws.Range("I48").Value = 8
$ws.Range("J48").Value = 9
$ws.Range("I49").Value = 8
$ws.Range("J49").Value = 8
$ws.Range("I50").Value = 9
$ws.Range("J50").Value = 9
$ws.Range("I51").Value = 7
$ws.Range("J51").Value = 9
$ws.Range("I52").Value = 9
$ws.Range("J52").Value = 9
$ws.Range("I53").Value = 9
$ws.Range("J53").Value = 9
$ws.Range("I54").Value = 8
$ws.Range("J54").Value = 8
$ws.Range("I55").Value = 9
$ws.Range("J55").Value = 9
$ws.Range("I56").Value = 9
$ws.Range("J56").Value = 9
$ws.Range("I57").Value = 9
$ws.Range("J57").Value = 9
$ws.Range("I58").Value = 8
$ws.Range("J58").Value = 9
$ws.Range("I59").Value = 9
$ws.Range("J59").Value = 9
$ws.Range("I60").Value = 10
$ws.Range("J60").Value = 10
$ws.Range("I61").Value = 9
$ws.Range("J61").Value = 10
$ws.Range("I62").Value = 6
$ws.Range("J62").Value = 7
$ws.Range("I63").Value = 5
$ws.Range("J63").Value = 7
$ws.Range("I64").Value = 7
$ws.Range("J64").Value = 8
$ws.Range("I65").Value = 8
$ws.Range("J65").Value = 8
$ws.Range("I66").Value = 7
$ws.Range("J66").Value = 7
$ws.Range("I67").Value = 9
$ws.Range("J67").Value = 9
$ws.Range("I68").Value = 8
$ws.Range("J68").Value = 8
$ws.Range("I69").Value = 9
$ws.Range("J69").Value = 9
$ws.Range("I70").Value = 5
$ws.Range("J70").Value = 5
$ws.Range("I71").Value = 6
$ws.Range("J71").Value = 7
$ws.Range("I72").Value = 3
$ws.Range("J72").Value = 5
$ws.Range("I73").Value = 7
$ws.Range("J73").Value = 7
$ws.Range("I74").Value = 6
$ws.Range("J74").Value = 6
$ws.Range("I75").Value = 6
$ws.Range("J75").Value = 6
$ws.Range("I76").Value = 5
$ws.Range("J76").Value = 6
$ws.Range("I77").Value = 7
$ws.Range("J77").Value = 7
$ws.Range("I78").Value = 6
$ws.Range("J78").Value = 6
$ws.Range("I79").Value = 8
$ws.Range("J79").Value = 8
$ws.Range("I80").Value = 8
$ws.Range("J80").Value = 8
$ws.Range("I81").Value = 7
$ws.Range("J81").Value = 7
$ws.Range("I82").Value = 4
$ws.Range("J82").Value = 4
